# Re-run of the NATMI TPM pipeline for Ifnb1-Ifnar1: the sending-cluster set now
# includes "ECs" (in addition to Inflammatory-Mac / Resolving-Mac), and every
# expression/specificity metric is recomputed against the new TPM values.
# Each of the 3 sending clusters is paired with the 5 target clusters
# (ECs, FAPs, Inflammatory-Mac, MuSCs, Resolving-Mac), for 15 data rows total.
$data = @(
    @("ECs", "Ifnb1", "Ifnar1", "ECs", 1, 0.3333333333333333, 0.002067, 0.006201, 0.004794555500401285, 0.004794555500401286, 3, 1, 16.82700666666667, 50.48102, 0.1252829202765328, 0.1283458604125, 0.03478142277999999, 0.31303280502, 0.000600675914518186, 0.0006153613509944874),
    @("ECs", "Ifnb1", "Ifnar1", "FAPs", 1, 0.3333333333333333, 0.002067, 0.006201, 0.004794555500401285, 0.004794555500401286, 3, 1, 20.117406, 60.35221800000001, 0.1497810883418348, 0.1534429642470134, 0.041582678202, 0.3742441038180001, 0.0007181337409654349, 0.0007356908082283959),
    @("ECs", "Ifnb1", "Ifnar1", "Inflammatory-Mac", 1, 0.3333333333333333, 0.002067, 0.006201, 0.004794555500401285, 0.004794555500401286, 3, 1, 51.48756166666666, 154.462685, 0.3833428138051528, 0.3927148502140002, 0.106424789965, 0.957823109685, 0.001837958396468801, 0.001882893145182802),
    @("ECs", "Ifnb1", "Ifnar1", "MuSCs", 1, 0.3333333333333333, 0.002067, 0.006201, 0.004794555500401285, 0.004794555500401286, 2, 1, 9.6159655, 19.231931, 0.07159420941484877, 0.04889637197482995, 0.0198762006885, 0.119257204131, 0.0003432624105468447, 0.0002344363692015883),
    @("ECs", "Ifnb1", "Ifnar1", "Resolving-Mac", 1, 0.3333333333333333, 0.002067, 0.006201, 0.004794555500401285, 0.004794555500401286, 3, 1, 36.26411666666667, 108.79235, 0.2699989681616309, 0.2765999531516566, 0.07495792914999999, 0.67462136235, 0.001294525037902019, 0.001326173826794013),
    @("Inflammatory-Mac", "Ifnb1", "Ifnar1", "ECs", 2, 0.6666666666666666, 0.3352766666666667, 1.00583, 0.7776983968664127, 0.7776983968664127, 3, 1, 16.82700666666667, 50.48102, 0.1252829202765328, 0.1283458604125, 5.641702705177777, 50.7753243466, 0.09743232625380215, 0.0998143698872416),
    @("Inflammatory-Mac", "Ifnb1", "Ifnar1", "FAPs", 2, 0.6666666666666666, 0.3352766666666667, 1.00583, 0.7776983968664127, 0.7776983968664127, 3, 1, 20.117406, 60.35221800000001, 0.1497810883418348, 0.1534429642470134, 6.744896825660001, 60.70407143094, 0.1164845122843515, 0.1193323473053326),
    @("Inflammatory-Mac", "Ifnb1", "Ifnar1", "Inflammatory-Mac", 2, 0.6666666666666666, 0.3352766666666667, 1.00583, 0.7776983968664127, 0.7776983968664127, 3, 1, 51.48756166666666, 154.462685, 0.3833428138051528, 0.3927148502140002, 17.26257805039445, 155.36320245355, 0.2981250917465271, 0.3054137094370614),
    @("Inflammatory-Mac", "Ifnb1", "Ifnar1", "MuSCs", 2, 0.6666666666666666, 0.3352766666666667, 1.00583, 0.7776983968664127, 0.7776983968664127, 2, 1, 9.6159655, 19.231931, 0.07159420941484877, 0.04889637197482995, 3.224008859621667, 19.34405315773, 0.05567870188684612, 0.03802663009740904),
    @("Inflammatory-Mac", "Ifnb1", "Ifnar1", "Resolving-Mac", 2, 0.6666666666666666, 0.3352766666666667, 1.00583, 0.7776983968664127, 0.7776983968664127, 3, 1, 36.26411666666667, 108.79235, 0.2699989681616309, 0.2765999531516566, 12.15851215561111, 109.4266094005, 0.2099777646948859, 0.2151113401393682),
    @("Resolving-Mac", "Ifnb1", "Ifnar1", "ECs", 1, 0.3333333333333333, 0.09377033333333333, 0.281311, 0.217507047633186, 0.217507047633186, 3, 1, 16.82700666666667, 50.48102, 0.1252829202765328, 0.1283458604125, 1.577874024135555, 14.20086621722, 0.02724991810821245, 0.02791612917426387),
    @("Resolving-Mac", "Ifnb1", "Ifnar1", "FAPs", 1, 0.3333333333333333, 0.09377033333333333, 0.281311, 0.217507047633186, 0.217507047633186, 3, 1, 20.117406, 60.35221800000001, 0.1497810883418348, 0.1534429642470134, 1.886415866422, 16.977742797798, 0.03257844231651789, 0.03337492613345239),
    @("Resolving-Mac", "Ifnb1", "Ifnar1", "Inflammatory-Mac", 1, 0.3333333333333333, 0.09377033333333333, 0.281311, 0.217507047633186, 0.217507047633186, 3, 1, 51.48756166666666, 154.462685, 0.3833428138051528, 0.3927148502140002, 4.828005820003889, 43.45205238003499, 0.08337976366215691, 0.08541824763175604),
    @("Resolving-Mac", "Ifnb1", "Ifnar1", "MuSCs", 1, 0.3333333333333333, 0.09377033333333333, 0.281311, 0.217507047633186, 0.217507047633186, 2, 1, 9.6159655, 19.231931, 0.07159420941484877, 0.04889637197482995, 0.9016922902568333, 5.410153741540999, 0.0155722451174558, 0.01063530550821932),
    @("Resolving-Mac", "Ifnb1", "Ifnar1", "Resolving-Mac", 1, 0.3333333333333333, 0.09377033333333333, 0.281311, 0.217507047633186, 0.217507047633186, 3, 1, 36.26411666666667, 108.79235, 0.2699989681616309, 0.2765999531516566, 3.400498307872222, 30.60448477085, 0.05872667842884291, 0.06016243918549437),
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowCount = $data.Count
$colCount = $data[0].Count

# Overwrite/extend rows 2..16 (A:T) with the refreshed TPM results.
for ($i = 0; $i -lt $rowCount; $i++) {
    $rowNum = $i + 2
    for ($j = 0; $j -lt $colCount; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $data[$i][$j]
    }
}
